$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the commit diff: a price/volume refresh plus a
# handful of coin rows that changed rank order (Coin/Link/Price/Volume move
# together as a block). NumberFormat is forced to text ("@") before each
# assignment so Excel does not reinterpret numeric-looking strings (prices
# like "1.002" or "103.53") as real numbers, then the style is reset back to
# "Normal" so no stray number-format style sticks to the cell.

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '26.824.70'
$c.Style = "Normal"
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  -1.90%  '
$c.Style = "Normal"
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '1.801.96'
$c.Style = "Normal"
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  -1.09%  '
$c.Style = "Normal"
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '
$c.Style = "Normal"
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '308.83'
$c.Style = "Normal"
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -1.67%  '
$c.Style = "Normal"
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '
$c.Style = "Normal"
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.4597'
$c.Style = "Normal"
$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +3.53%  '
$c.Style = "Normal"
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3721'
$c.Style = "Normal"
$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  -0.76%  '
$c.Style = "Normal"
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.07349'
$c.Style = "Normal"
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  -1.71%  '
$c.Style = "Normal"
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.8712'
$c.Style = "Normal"
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  -1.45%  '
$c.Style = "Normal"
$c = $ws.Range('B12')
$c.NumberFormat = "@"
$c.Value = 'WrappedEther'
$c.Style = "Normal"
$c = $ws.Range('C12')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c.Style = "Normal"
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '1.868.38'
$c.Style = "Normal"
$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  +2.49%  '
$c.Style = "Normal"
$c = $ws.Range('B13')
$c.NumberFormat = "@"
$c.Value = 'Polkadot'
$c.Style = "Normal"
$c = $ws.Range('C13')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c.Style = "Normal"
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '5.354'
$c.Style = "Normal"
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -0.85%  '
$c.Style = "Normal"
$c = $ws.Range('B14')
$c.NumberFormat = "@"
$c.Value = 'Chainlink'
$c.Style = "Normal"
$c = $ws.Range('C14')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c.Style = "Normal"
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '6.510'
$c.Style = "Normal"
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  -3.50%  '
$c.Style = "Normal"
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.07047'
$c.Style = "Normal"
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -1.11%  '
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '91.20'
$c.Style = "Normal"
$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  -2.49%  '
$c.Style = "Normal"
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.000008702'
$c.Style = "Normal"
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -0.59%  '
$c.Style = "Normal"
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.Style = "Normal"
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '14.65'
$c.Style = "Normal"
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  -3.23%  '
$c.Style = "Normal"
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '26.840.39'
$c.Style = "Normal"
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -1.86%  '
$c.Style = "Normal"
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '5.293'
$c.Style = "Normal"
$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -1.63%  '
$c.Style = "Normal"
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '10.68'
$c.Style = "Normal"
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  -2.21%  '
$c.Style = "Normal"
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.058.21'
$c.Style = "Normal"
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.Style = "Normal"
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '1.908'
$c.Style = "Normal"
$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  -2.68%  '
$c.Style = "Normal"
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '151.29'
$c.Style = "Normal"
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  +0.07%  '
$c.Style = "Normal"
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '18.36'
$c.Style = "Normal"
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  -1.28%  '
$c.Style = "Normal"
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.134'
$c.Style = "Normal"
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  -8.25%  '
$c.Style = "Normal"
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '5.250'
$c.Style = "Normal"
$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -1.93%  '
$c.Style = "Normal"
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '115.49'
$c.Style = "Normal"
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -1.80%  '
$c.Style = "Normal"
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '0.08881'
$c.Style = "Normal"
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  +0.10%  '
$c.Style = "Normal"
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.7594'
$c.Style = "Normal"
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  -3.14%  '
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '1.151'
$c.Style = "Normal"
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  -4.20%  '
$c.Style = "Normal"
$c = $ws.Range('B34')
$c.NumberFormat = "@"
$c.Value = 'HuobiToken'
$c.Style = "Normal"
$c = $ws.Range('C34')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c.Style = "Normal"
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '2.922'
$c.Style = "Normal"
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  +0.71%  '
$c.Style = "Normal"
$c = $ws.Range('B35')
$c.NumberFormat = "@"
$c.Value = 'Filecoin'
$c.Style = "Normal"
$c = $ws.Range('C35')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.Style = "Normal"
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '4.451'
$c.Style = "Normal"
$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  -3.35%  '
$c.Style = "Normal"
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '1.105'
$c.Style = "Normal"
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.Style = "Normal"
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.01947'
$c.Style = "Normal"
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -2.36%  '
$c.Style = "Normal"
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.05233'
$c.Style = "Normal"
$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  -1.35%  '
$c.Style = "Normal"
$c = $ws.Range('B40')
$c.NumberFormat = "@"
$c.Value = 'FraxShare'
$c.Style = "Normal"
$c = $ws.Range('C40')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c.Style = "Normal"
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '7.259'
$c.Style = "Normal"
$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '
$c.Style = "Normal"
$c = $ws.Range('B41')
$c.NumberFormat = "@"
$c.Value = 'MXToken'
$c.Style = "Normal"
$c = $ws.Range('C41')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c.Style = "Normal"
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '2.908'
$c.Style = "Normal"
$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  +1.78%  '
$c.Style = "Normal"
$c = $ws.Range('B42')
$c.NumberFormat = "@"
$c.Value = 'TheSandbox'
$c.Style = "Normal"
$c = $ws.Range('C42')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c.Style = "Normal"
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.5276'
$c.Style = "Normal"
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  -0.71%  '
$c.Style = "Normal"
$c = $ws.Range('B43')
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c.Style = "Normal"
$c = $ws.Range('C43')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c.Style = "Normal"
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.337'
$c.Style = "Normal"
$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  +0.75%  '
$c.Style = "Normal"
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.1653'
$c.Style = "Normal"
$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  -3.45%  '
$c.Style = "Normal"
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '8.491'
$c.Style = "Normal"
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  -1.67%  '
$c.Style = "Normal"
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.5013'
$c.Style = "Normal"
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -1.26%  '
$c.Style = "Normal"
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '10.27'
$c.Style = "Normal"
$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -3.23%  '
$c.Style = "Normal"
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  +0.04%  '
$c.Style = "Normal"
$c = $ws.Range('B49')
$c.NumberFormat = "@"
$c.Value = 'Quant'
$c.Style = "Normal"
$c = $ws.Range('C49')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.Style = "Normal"
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '103.53'
$c.Style = "Normal"
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -1.67%  '
$c.Style = "Normal"
$c = $ws.Range('B50')
$c.NumberFormat = "@"
$c.Value = 'NEARProtocol'
$c.Style = "Normal"
$c = $ws.Range('C50')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c.Style = "Normal"
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '1.658'
$c.Style = "Normal"
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  -1.96%  '
$c.Style = "Normal"
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.06294'
$c.Style = "Normal"
$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -1.57%  '
$c.Style = "Normal"
